# Fix bugs on medical records: append the missing AP004 appointment row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Appointment ID, B=Appointment Date & Time, C=Patient ID,
#          D=Doctor ID, E=Appointment Status, F=Outcome Record ID
# Fill order chosen to match how the unique shared-string list ends up
# ordered in the saved file (AP004, OR004, date, then status).
$ws.Range("A5").Value = "AP004"
$ws.Range("F5").Value = "OR004"
$ws.Range("B5").Value = "2024-11-19T19:10:12.981009"

# Status was typed as "Pending" first and then corrected to "Available".
$ws.Range("E5").Value = "Pending"
$ws.Range("E5").Value = "Available"

$ws.Range("C5").Value = "N/A"
$ws.Range("D5").Value = "H002"

# Leave the selection where it ended up after entering the new row.
$ws.Range("D6").Select()
